$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert the new "active learning" concentration rows.
# We work from the bottom of the sheet upward so that row numbers used below
# always refer to the not-yet-shifted original layout. Where several rows are
# inserted at the very same point, they are written in the reverse of their
# final order (each subsequent insert pushes the earlier one down).
# ---------------------------------------------------------------------------

# New SecYE concentration (0.625) between SecYE=0 and SecYE=1.25 (orig row 37)
$ws.Range("A37").EntireRow.Insert()
$ws.Range("A37").Value = "SecYE"
$ws.Range("B37").Value = "SecYE"
$ws.Range("C37").Value = 0.625

# Two new "Liposome_name" rows + new K concentration (235), all inserted
# right before the old SecYE=0 row (orig row 36) in reverse final order.
$ws.Range("A36").EntireRow.Insert()
$ws.Range("A36").Value = "Liposome_name"
$ws.Range("B36").Value = "Liposome_name"
$ws.Range("C36").Value = "DMPC"

$ws.Range("A36").EntireRow.Insert()
$ws.Range("A36").Value = "Liposome_name"
$ws.Range("B36").Value = "Liposome_name"
$ws.Range("C36").Value = "DOPC"

$ws.Range("A36").EntireRow.Insert()
$ws.Range("A36").Value = "K"
$ws.Range("B36").Value = "K"
$ws.Range("C36").Value = 235

# New Mg concentration (23) between Mg=20 and K=85 (orig row 33)
$ws.Range("A33").EntireRow.Insert()
$ws.Range("A33").Value = "Mg"
$ws.Range("B33").Value = "Mg"
$ws.Range("C33").Value = 23

# New Mg concentration (17) between Mg=14 and Mg=20 (orig row 32)
$ws.Range("A32").EntireRow.Insert()
$ws.Range("A32").Value = "Mg"
$ws.Range("B32").Value = "Mg"
$ws.Range("C32").Value = 17

# New Mg concentration (11) between Mg=8 and Mg=14 (orig row 31)
$ws.Range("A31").EntireRow.Insert()
$ws.Range("A31").Value = "Mg"
$ws.Range("B31").Value = "Mg"
$ws.Range("C31").Value = 11

# ---------------------------------------------------------------------------
# Remove the old "Liposome" rows (they're superseded by the Liposome_name
# rows above). After all the inserts above they now sit at rows 48 and 49.
# ---------------------------------------------------------------------------
$ws.Range("A49").EntireRow.Delete()
$ws.Range("A48").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Refresh the AutoFilter so it covers the new data range. The AutoFilter
# always snaps to the sheet's used range, so temporarily clear the very last
# row (PEG/PEG/2, row 47) before (re)applying the filter -- this reproduces
# the source file, where the filter only grew to row 46 -- then restore it.
# ---------------------------------------------------------------------------
$ws.Range("A47:C47").ClearContents()

$ws.AutoFilterMode = $false
$ws.Range("A1:C46").AutoFilter()

$ws.Range("A47").Value = "PEG"
$ws.Range("B47").Value = "PEG"
$ws.Range("C47").Value = 2

$ws.Range("C42").Select()

# Keep the "_xlnm._FilterDatabase" defined name in sync with the new range.
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.Name() -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$46"
    }
}
